$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.854.05"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").Value = "'3.212.21"
$ws.Range("E3").Value = "  -3.55%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'581.43"
$ws.Range("E5").Value = "  -3.72%  "

$ws.Range("D6").Value = "'140.44"
$ws.Range("E6").Value = "  -14.21%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'3.206.85"
$ws.Range("E8").Value = "  -3.76%  "

$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  -9.08%  "

$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  -12.99%  "

$ws.Range("D11").Value = "'6.49"
$ws.Range("E11").Value = "  -2.77%  "

$ws.Range("D12").Value = "'0.477"
$ws.Range("E12").Value = "  -10.01%  "

$ws.Range("D13").Value = "'0.0000232"
$ws.Range("E13").Value = "  -10.02%  "

$ws.Range("D14").Value = "'35.81"
$ws.Range("E14").Value = "  -14.14%  "

$ws.Range("D15").Value = "'3.733.13"
$ws.Range("E15").Value = "  -3.67%  "

$ws.Range("D16").Value = "'66.911.21"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").Value = "'3.210.36"
$ws.Range("E17").Value = "  -3.74%  "

$ws.Range("E18").Value = "  -4.99%  "

$ws.Range("D19").Value = "'6.76"
$ws.Range("E19").Value = "  -12.40%  "

$ws.Range("D20").Value = "'497.48"
$ws.Range("E20").Value = "  -10.81%  "

$ws.Range("D21").Value = "'14.24"
$ws.Range("E21").Value = "  -11.78%  "

$ws.Range("D22").Value = "'0.714"
$ws.Range("E22").Value = "  -11.04%  "

$ws.Range("D23").Value = "'7.30"
$ws.Range("E23").Value = "  -13.43%  "

$ws.Range("D24").Value = "'81.99"
$ws.Range("E24").Value = "  -8.65%  "

$ws.Range("D25").Value = "'12.83"
$ws.Range("E25").Value = "  -10.06%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").Value = "'3.08"
$ws.Range("E27").Value = "  -11.34%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.02"
$ws.Range("E28").Value = "  -11.13%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'27.64"
$ws.Range("E29").Value = "  -11.15%  "

$ws.Range("D30").Value = "'7.50"
$ws.Range("E30").Value = "  -8.35%  "

$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").Value = "'2.50"
$ws.Range("E32").Value = "  -5.74%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'6.07"
$ws.Range("E34").Value = "  -17.13%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'495.64"
$ws.Range("E35").Value = "  -13.32%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'54.30"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("D37").Value = "'5.28"
$ws.Range("E37").Value = "  -15.26%  "

$ws.Range("D38").Value = "'0.0417"
$ws.Range("E38").Value = "  -6.97%  "

$ws.Range("D39").Value = "'0.0813"
$ws.Range("E39").Value = "  -10.30%  "

$ws.Range("D40").Value = "'8.45"
$ws.Range("E40").Value = "  -14.74%  "

$ws.Range("E41").Value = "  -14.39%  "

$ws.Range("D42").Value = "'2.826.30"
$ws.Range("E42").Value = "  -7.96%  "

$ws.Range("D43").Value = "'2.53"
$ws.Range("E43").Value = "  -13.34%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.250"
$ws.Range("E44").Value = "  -10.14%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'25.06"
$ws.Range("E46").Value = "  -13.21%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.03"
$ws.Range("E47").Value = "  -10.12%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'121.46"
$ws.Range("E48").Value = "  -6.18%  "

$ws.Range("D49").Value = "'0.0₃0526"
$ws.Range("E49").Value = "  -15.83%  "

$ws.Range("D50").Value = "'0.109"
$ws.Range("E50").Value = "  -9.60%  "

$ws.Range("D51").Value = "'2.12"
$ws.Range("E51").Value = "  -20.46%  "
